$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing "sum" header (G1) into the new H1 header cell,
# then set its text. This reuses the same header style (bold, bordered, centered)
# rather than creating a new style entry.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New "Save" data column value for row 2 (plain numeric cell, no special style)
$ws.Range("H2").Value = 0
